$d = $word.ActiveDocument

# Locate the first paragraph (ends with "Заявление на заселение.")
$p1 = $d.Paragraphs(1)
$r = $p1.Range

# Insert a new paragraph right after it.
$r.InsertParagraphAfter()

# The newly created paragraph is now paragraph #2; grab its range.
$p2 = $d.Paragraphs(2)
$newRange = $p2.Range

# Apply bold paragraph mark formatting matching the first paragraph, then set the text.
$newRange.Text = "Действующее лицо: гражданин."
$newRange.Font.Bold = $true

$newRange.ParagraphFormat.SpaceAfter = 0
$newRange.ParagraphFormat.LineSpacingRule = 1
$newRange.ParagraphFormat.LineSpacing = 13.8
$newRange.ParagraphFormat.FirstLineIndent = 0
$newRange.ParagraphFormat.Alignment = 0
